# chore: update Sheets via scheduled runner
# Refreshes the market-price derived columns (H:N) for the affected Leve
# rows across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H42").Value = 91
$ws_ALC.Range("I42").Value = 79.59999999999999
$ws_ALC.Range("J42").Value = 97.333336
$ws_ALC.Range("K42").Value = 238.8
$ws_ALC.Range("L42").Value = 292.000008
$ws_ALC.Range("M42").Value = -8.799999999999983
$ws_ALC.Range("N42").Value = -752.000008

$ws_ALC.Range("H100").Value = 50001268
$ws_ALC.Range("I100").Value = 1633.3334
$ws_ALC.Range("J100").Value = 125000720
$ws_ALC.Range("K100").Value = 1633.3334
$ws_ALC.Range("L100").Value = 125000720
$ws_ALC.Range("M100").Value = -1092.3334
$ws_ALC.Range("N100").Value = -125001802

$ws_ALC.Range("H137").Value = 2415.5
$ws_ALC.Range("I137").Value = 1979.2
$ws_ALC.Range("J137").Value = 2727.1428
$ws_ALC.Range("K137").Value = 5937.6
$ws_ALC.Range("L137").Value = 8181.428400000001
$ws_ALC.Range("M137").Value = -3387.6
$ws_ALC.Range("N137").Value = -13281.4284

$ws_ALC.Range("H138").Value = 1872.3715
$ws_ALC.Range("I138").Value = 1387
$ws_ALC.Range("J138").Value = 2281.1052
$ws_ALC.Range("K138").Value = 4161
$ws_ALC.Range("L138").Value = 6843.3156
$ws_ALC.Range("M138").Value = 979
$ws_ALC.Range("N138").Value = -17123.3156

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H52").Value = 16244.444
$ws_ARM.Range("J52").Value = 16244.444
$ws_ARM.Range("L52").Value = 16244.444
$ws_ARM.Range("N52").Value = -16880.444

$ws_ARM.Range("H55").Value = 14550.75
$ws_ARM.Range("J55").Value = 14550.75
$ws_ARM.Range("L55").Value = 14550.75
$ws_ARM.Range("N55").Value = -15180.75

$ws_ARM.Range("H61").Value = 1741.5714
$ws_ARM.Range("I61").Value = 1904.5714
$ws_ARM.Range("J61").Value = 1578.5714
$ws_ARM.Range("K61").Value = 1904.5714
$ws_ARM.Range("L61").Value = 1578.5714
$ws_ARM.Range("M61").Value = -1692.5714
$ws_ARM.Range("N61").Value = -2002.5714

$ws_ARM.Range("H74").Value = 1075.2903
$ws_ARM.Range("I74").Value = 803.5
$ws_ARM.Range("J74").Value = 2007.1428
$ws_ARM.Range("K74").Value = 803.5
$ws_ARM.Range("L74").Value = 2007.1428
$ws_ARM.Range("M74").Value = 70.5
$ws_ARM.Range("N74").Value = -3755.1428

$ws_ARM.Range("H77").Value = 1075.2903
$ws_ARM.Range("I77").Value = 803.5
$ws_ARM.Range("J77").Value = 2007.1428
$ws_ARM.Range("K77").Value = 4017.5
$ws_ARM.Range("L77").Value = 10035.714
$ws_ARM.Range("M77").Value = 350.5
$ws_ARM.Range("N77").Value = -18771.714

$ws_ARM.Range("H80").Value = 22531.428
$ws_ARM.Range("J80").Value = 22531.428
$ws_ARM.Range("L80").Value = 22531.428
$ws_ARM.Range("N80").Value = -24527.428

$ws_ARM.Range("H83").Value = 22531.428
$ws_ARM.Range("J83").Value = 22531.428
$ws_ARM.Range("L83").Value = 67594.284
$ws_ARM.Range("N83").Value = -77578.284

$ws_ARM.Range("H110").Value = 7568
$ws_ARM.Range("I110").Value = 8296.333000000001
$ws_ARM.Range("J110").Value = 1013
$ws_ARM.Range("K110").Value = 8296.333000000001
$ws_ARM.Range("L110").Value = 1013
$ws_ARM.Range("M110").Value = -6251.333000000001
$ws_ARM.Range("N110").Value = -5103

$ws_ARM.Range("H136").Value = 1741.5714
$ws_ARM.Range("I136").Value = 1904.5714
$ws_ARM.Range("J136").Value = 1578.5714
$ws_ARM.Range("K136").Value = 5713.7142
$ws_ARM.Range("L136").Value = 4735.7142
$ws_ARM.Range("M136").Value = -3163.7142
$ws_ARM.Range("N136").Value = -9835.7142

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H59").Value = 53000
$ws_BSM.Range("J59").Value = 53000
$ws_BSM.Range("L59").Value = 53000
$ws_BSM.Range("N59").Value = -54694

$ws_BSM.Range("H82").Value = 23641.39
$ws_BSM.Range("I82").Value = 7338.6665
$ws_BSM.Range("J82").Value = 29395.295
$ws_BSM.Range("K82").Value = 7338.6665
$ws_BSM.Range("L82").Value = 29395.295
$ws_BSM.Range("M82").Value = -6955.6665
$ws_BSM.Range("N82").Value = -30161.295

$ws_BSM.Range("H85").Value = 23641.39
$ws_BSM.Range("I85").Value = 7338.6665
$ws_BSM.Range("J85").Value = 29395.295
$ws_BSM.Range("K85").Value = 7338.6665
$ws_BSM.Range("L85").Value = 29395.295
$ws_BSM.Range("M85").Value = -6012.6665
$ws_BSM.Range("N85").Value = -32047.295

$ws_BSM.Range("H86").Value = 2117.25
$ws_BSM.Range("I86").Value = 2191.238
$ws_BSM.Range("J86").Value = 1895.2858
$ws_BSM.Range("K86").Value = 2191.238
$ws_BSM.Range("L86").Value = 1895.2858
$ws_BSM.Range("M86").Value = -1068.238
$ws_BSM.Range("N86").Value = -4141.2858

$ws_BSM.Range("H89").Value = 2117.25
$ws_BSM.Range("I89").Value = 2191.238
$ws_BSM.Range("J89").Value = 1895.2858
$ws_BSM.Range("K89").Value = 10956.19
$ws_BSM.Range("L89").Value = 9476.429
$ws_BSM.Range("M89").Value = -5340.189999999999
$ws_BSM.Range("N89").Value = -20708.429

$ws_BSM.Range("H107").Value = 4901.4287
$ws_BSM.Range("I107").Value = 4785.7896
$ws_BSM.Range("J107").Value = 6000
$ws_BSM.Range("K107").Value = 4785.7896
$ws_BSM.Range("L107").Value = 6000
$ws_BSM.Range("M107").Value = -2865.7896
$ws_BSM.Range("N107").Value = -9840

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H16").Value = 2300.48
$ws_CRP.Range("I16").Value = 2970.353
$ws_CRP.Range("J16").Value = 877
$ws_CRP.Range("K16").Value = 2970.353
$ws_CRP.Range("L16").Value = 877
$ws_CRP.Range("M16").Value = -2683.353
$ws_CRP.Range("N16").Value = -1451

$ws_CRP.Range("H31").Value = 2724.5938
$ws_CRP.Range("I31").Value = 1915.15
$ws_CRP.Range("J31").Value = 4073.6667
$ws_CRP.Range("K31").Value = 1915.15
$ws_CRP.Range("L31").Value = 4073.6667
$ws_CRP.Range("M31").Value = -1620.15
$ws_CRP.Range("N31").Value = -4663.6667

$ws_CRP.Range("H34").Value = 2724.5938
$ws_CRP.Range("I34").Value = 1915.15
$ws_CRP.Range("J34").Value = 4073.6667
$ws_CRP.Range("K34").Value = 1915.15
$ws_CRP.Range("L34").Value = 4073.6667
$ws_CRP.Range("M34").Value = -1713.15
$ws_CRP.Range("N34").Value = -4477.6667

$ws_CRP.Range("H41").Value = 16138.75
$ws_CRP.Range("J41").Value = 21351.666
$ws_CRP.Range("L41").Value = 21351.666
$ws_CRP.Range("N41").Value = -22207.666

$ws_CRP.Range("H50").Value = 8989.429
$ws_CRP.Range("J50").Value = 8989.429
$ws_CRP.Range("L50").Value = 8989.429
$ws_CRP.Range("N50").Value = -10239.429

$ws_CRP.Range("H51").Value = 9382.4
$ws_CRP.Range("J51").Value = 9382.4
$ws_CRP.Range("L51").Value = 9382.4
$ws_CRP.Range("N51").Value = -10854.4

$ws_CRP.Range("H60").Value = 32631.5
$ws_CRP.Range("I60").Value = 3093
$ws_CRP.Range("J60").Value = 36851.285
$ws_CRP.Range("K60").Value = 3093
$ws_CRP.Range("L60").Value = 36851.285
$ws_CRP.Range("M60").Value = -2582
$ws_CRP.Range("N60").Value = -37873.285

$ws_CRP.Range("H61").Value = 9382.4
$ws_CRP.Range("J61").Value = 9382.4
$ws_CRP.Range("L61").Value = 9382.4
$ws_CRP.Range("N61").Value = -10078.4

$ws_CRP.Range("H68").Value = 16794.75
$ws_CRP.Range("J68").Value = 16794.75
$ws_CRP.Range("L68").Value = 16794.75
$ws_CRP.Range("N68").Value = -18292.75

$ws_CRP.Range("H71").Value = 16794.75
$ws_CRP.Range("J71").Value = 16794.75
$ws_CRP.Range("L71").Value = 50384.25
$ws_CRP.Range("N71").Value = -57872.25

$ws_CRP.Range("H103").Value = 12188.889
$ws_CRP.Range("I103").Value = 9087.5
$ws_CRP.Range("J103").Value = 37000
$ws_CRP.Range("K103").Value = 9087.5
$ws_CRP.Range("L103").Value = 37000
$ws_CRP.Range("M103").Value = -7915.5
$ws_CRP.Range("N103").Value = -39344

$ws_CRP.Range("H113").Value = 2300.48
$ws_CRP.Range("I113").Value = 2970.353
$ws_CRP.Range("J113").Value = 877
$ws_CRP.Range("K113").Value = 2970.353
$ws_CRP.Range("L113").Value = 877
$ws_CRP.Range("M113").Value = -800.3530000000001
$ws_CRP.Range("N113").Value = -5217

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 321.56522
$ws_CUL.Range("I5").Value = 297.42856
$ws_CUL.Range("J5").Value = 575
$ws_CUL.Range("K5").Value = 892.28568
$ws_CUL.Range("L5").Value = 1725
$ws_CUL.Range("M5").Value = -780.28568
$ws_CUL.Range("N5").Value = -1949

$ws_CUL.Range("H69").Value = 3800
$ws_CUL.Range("I69").Value = 500
$ws_CUL.Range("J69").Value = 4460
$ws_CUL.Range("K69").Value = 1500
$ws_CUL.Range("L69").Value = 13380
$ws_CUL.Range("M69").Value = -689
$ws_CUL.Range("N69").Value = -15002

$ws_CUL.Range("H72").Value = 3800
$ws_CUL.Range("I72").Value = 500
$ws_CUL.Range("J72").Value = 4460
$ws_CUL.Range("K72").Value = 4500
$ws_CUL.Range("L72").Value = 40140
$ws_CUL.Range("M72").Value = -444
$ws_CUL.Range("N72").Value = -48252

$ws_CUL.Range("H115").Value = 1961.1428
$ws_CUL.Range("I115").Value = 1382
$ws_CUL.Range("J115").Value = 2733.3333
$ws_CUL.Range("K115").Value = 4146
$ws_CUL.Range("L115").Value = 8199.999899999999
$ws_CUL.Range("M115").Value = -2971
$ws_CUL.Range("N115").Value = -10549.9999

$ws_CUL.Range("H135").Value = 321.56522
$ws_CUL.Range("I135").Value = 297.42856
$ws_CUL.Range("J135").Value = 575
$ws_CUL.Range("K135").Value = 2676.85704
$ws_CUL.Range("L135").Value = 5175
$ws_CUL.Range("M135").Value = -141.8570399999999
$ws_CUL.Range("N135").Value = -10245

$ws_CUL.Range("H141").Value = 9902.223
$ws_CUL.Range("I141").Value = 4224
$ws_CUL.Range("J141").Value = 17000
$ws_CUL.Range("K141").Value = 12672
$ws_CUL.Range("L141").Value = 51000
$ws_CUL.Range("M141").Value = -7492
$ws_CUL.Range("N141").Value = -61360

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H57").Value = 4999
$ws_GSM.Range("J57").Value = 0
$ws_GSM.Range("L57").Value = 0
$ws_GSM.Range("N57").ClearContents()

$ws_GSM.Range("H113").Value = 1478.4
$ws_GSM.Range("I113").Value = 1463.4615
$ws_GSM.Range("J113").Value = 1506.1428
$ws_GSM.Range("K113").Value = 1463.4615
$ws_GSM.Range("L113").Value = 1506.1428
$ws_GSM.Range("M113").Value = 706.5385000000001
$ws_GSM.Range("N113").Value = -5846.1428

$ws_GSM.Range("H123").Value = 28073.2
$ws_GSM.Range("J123").Value = 28073.2
$ws_GSM.Range("L123").Value = 28073.2
$ws_GSM.Range("N123").Value = -32973.2

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H54").Value = 6038.5

$ws_WVR.Range("H96").Value = 1687.25
$ws_WVR.Range("I96").Value = 1700
$ws_WVR.Range("J96").Value = 1683
$ws_WVR.Range("K96").Value = 1700
$ws_WVR.Range("L96").Value = 1683
$ws_WVR.Range("M96").Value = -327
$ws_WVR.Range("N96").Value = -4429

$ws_WVR.Range("H109").Value = 0
$ws_WVR.Range("J109").Value = 0
$ws_WVR.Range("L109").Value = 0
$ws_WVR.Range("N109").ClearContents()
